$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) columns.
# NumberFormat is forced to Text ("@") before assignment so that
# numeric-looking strings (e.g. "1.00", "211.55") are preserved
# exactly as text instead of being auto-converted to numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.945.31'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.634.13'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.55'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.78%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.38'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.50%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.77%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.867.24'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.634.57'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.05%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.61%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.24'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.949.96'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '229.70'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.81'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.40%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.53%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.41%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.12'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.58%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '155.96'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.65%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.38%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.54'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.82%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.51%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.48%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.49%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.59%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.401.61'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.48%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.02'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.45%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.50%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.51%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.54%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.68%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.19%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '66.08'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.45%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.775.93'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.90%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.48'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.64%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.88%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.64'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.52%  '
